# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.809.79"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "1.628.68"
$ws.Range("E3").Value = "  -0.63%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5100"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2587"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06399"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.31%  "

# Row 12
$ws.Range("E12").Value = "  -0.22%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.626.64"

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.853.58"
$ws.Range("E14").Value = "  -0.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5591"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.45%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.48"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.35%  "

# Row 17
$ws.Range("D17").Value = "0.0₅7536"
$ws.Range("E17").Value = "  -2.73%  "

# Row 18
$ws.Range("D18").Value = "25.835.72"
$ws.Range("E18").Value = "  -0.59%  "

# Row 19
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.326"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.807"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.820"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1284"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "140.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.744"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.45"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.238"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04891"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.311"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.188"
$ws.Range("D33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.560"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.381"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8956"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.99%  "

# Row 37
$ws.Range("D37").Value = "1.132.05"
$ws.Range("E37").Value = "  +2.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.545"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5487"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01562"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.46%  "

# Row 41
$ws.Range("E41").Value = "  +0.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.593"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7966"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.02%  "

# Row 45
$ws.Range("D45").Value = "1.776.79"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("E46").Value = "  -7.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4435"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.99"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("E49").Value = "  -2.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.567"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.96%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.45%  "
